$wb = $excel.ActiveWorkbook

# The "想去人数" (interested-count) column F changed for the first five
# data rows. These two sheets ("展览" and "全部类型") carry identical
# data, so the same update is applied to both.
$sheetNames = @("展览", "全部类型")

$updates = @{
    2 = 335
    3 = 1420
    4 = 94
    5 = 75
    6 = 10
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
